$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 321.7
$ws.Range("C3").Value = 299.8
$ws.Range("C4").Value = 293.7
$ws.Range("C5").Value = 277.5
$ws.Range("C6").Value = 309.7
$ws.Range("C7").Value = 312
$ws.Range("C9").Value = 339.9
$ws.Range("C10").Value = 379.8
$ws.Range("C13").Value = 467.6
$ws.Range("C16").Value = 405.6
$ws.Range("C17").Value = 336.7
$ws.Range("C19").Value = 240.4
$ws.Range("C20").Value = 226.9
$ws.Range("C21").Value = 210.3
$ws.Range("C24").Value = 188.1
